$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the updated cells keep their original text-cell representation
# (these columns hold free-form strings like "331.71" or "0.45%", not numbers).
$cells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7","E8","D9","E9","D10","E10","E11","D12","E12","D13","E13","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","E23","D24","E24","D25","E25","E26","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","D46","E46","E47","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $cells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "331.71"
$ws.Range("E2").Value = "0.45%"
$ws.Range("D3").Value = "45.54"
$ws.Range("E3").Value = "2.63%"
$ws.Range("D4").Value = "5.564"
$ws.Range("E4").Value = "1.19%"
$ws.Range("D5").Value = "0.08365"
$ws.Range("E5").Value = "4.27%"
$ws.Range("D6").Value = "2.094"
$ws.Range("E6").Value = "-0.27%"
$ws.Range("D7").Value = "0.9905"
$ws.Range("E7").Value = "3.87%"
$ws.Range("E8").Value = "-0.76%"
$ws.Range("D9").Value = "0.1197"
$ws.Range("E9").Value = "4.43%"
$ws.Range("D10").Value = "0.1933"
$ws.Range("E10").Value = "1.37%"
$ws.Range("E11").Value = "0.68%"
$ws.Range("D12").Value = "0.09943"
$ws.Range("E12").Value = "0.08%"
$ws.Range("D13").Value = "0.04668"
$ws.Range("E13").Value = "-3.56%"
$ws.Range("E14").Value = "-0.55%"
$ws.Range("D15").Value = "0.001297"
$ws.Range("E15").Value = "1.89%"
$ws.Range("D16").Value = "0.005923"
$ws.Range("E16").Value = "-0.13%"
$ws.Range("D17").Value = "3.395"
$ws.Range("E17").Value = "0.79%"
$ws.Range("D18").Value = "4.441"
$ws.Range("E18").Value = "0.85%"
$ws.Range("D19").Value = "0.3368"
$ws.Range("E19").Value = "-2.57%"
$ws.Range("D20").Value = "0.1362"
$ws.Range("E20").Value = "-1.49%"
$ws.Range("D21").Value = "0.2563"
$ws.Range("E21").Value = "-0.76%"
$ws.Range("D22").Value = "0.04150"
$ws.Range("E22").Value = "1.87%"
$ws.Range("E23").Value = "1.43%"
$ws.Range("D24").Value = "0.004534"
$ws.Range("E24").Value = "4.15%"
$ws.Range("D25").Value = "0.0001302"
$ws.Range("E25").Value = "8.50%"
$ws.Range("E26").Value = "0.00%"
$ws.Range("D38").Value = "0.02699"
$ws.Range("E38").Value = "4.20%"
$ws.Range("D39").Value = "0.05744"
$ws.Range("E39").Value = "-0.94%"
$ws.Range("D40").Value = "0.007872"
$ws.Range("E40").Value = "4.11%"
$ws.Range("D41").Value = "0.1434"
$ws.Range("E41").Value = "2.23%"
$ws.Range("D42").Value = "0.007868"
$ws.Range("E42").Value = "9.36%"
$ws.Range("D43").Value = "0.002024"
$ws.Range("E43").Value = "0.39%"
$ws.Range("D44").Value = "0.008938"
$ws.Range("E44").Value = "-1.58%"
$ws.Range("D45").Value = "0.3408"
$ws.Range("D46").Value = "0.00007035"
$ws.Range("E46").Value = "0.44%"
$ws.Range("E47").Value = "0.15%"
$ws.Range("E48").Value = "0.26%"
$ws.Range("D49").Value = "0.003533"
$ws.Range("E49").Value = "0.05%"
$ws.Range("D50").Value = "0.003491"
$ws.Range("E50").Value = "-1.30%"
$ws.Range("D51").Value = "0.00002104"
$ws.Range("E51").Value = "0.15%"
